$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value = 7
$ws.Cells.Item(65, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(65, 3).Value = "Ñuble"
$ws.Cells.Item(65, 4).Value = 44413
$ws.Cells.Item(65, 5).Value = 16
$ws.Cells.Item(65, 6).Value = 100114013
$ws.Cells.Item(65, 7).Value = "Zanahoria"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 120
$ws.Cells.Item(65, 11).Value = 5000
$ws.Cells.Item(65, 12).Value = 5500
$ws.Cells.Item(65, 13).Value = 5250
$ws.Cells.Item(65, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(65, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(65, 16).Value = 262
$ws.Cells.Item(65, 17).Value = 20
$ws.Cells.Item(65, 18).Value = "Hortaliza"
